# ------------------------------------------------------------------
# Applies the SOC.pptx homework edit:
#   1. Move shape id=5 ("內容版面配置區 2") on slide 2 to a new position.
#   2. Add three new textboxes on slide 2 (after "文字方塊 11"):
#        - "文字方塊 9"  : #define XSCUGIC_INTR_PRIO_MASK ...
#        - "文字方塊 12" : #define XSCUGIC_PENDING_SET_OFFSET / CLR_OFFSET
#        - "文字方塊 14" : "pending暫存器的addresses" headline
#   3. Clear the title + body text on slide 5 ("Demo" / google-drive link).
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ====================================================================
# Slide 2
# ====================================================================
$slide2 = $p.Slides.Item(2)

# --- 1. Move the existing content placeholder (id=5) ----------------
$moved = $slide2.Shapes.Item(3)
$moved.Left = 12.3367
$moved.Top  = 154.6031

# --- 2a. New textbox "文字方塊 9" ------------------------------------
$tb9 = $slide2.Shapes.AddTextbox(1, 12.3367, 108.5973, 480.94906, 29.08126)
$tb9.Name = "文字方塊 9"
$tb9.Fill.Visible = 0
$tb9.TextFrame.WordWrap = -1
$tb9.TextFrame.AutoSize = 1

$tr9 = $tb9.TextFrame.TextRange
$tr9.Text = "#define XSCUGIC_INTR_PRIO_MASK`t0x000000F8U"

$r = $tr9.Characters(1, 7)                     # "#define"
$r.Font.Bold = $true
$r.Font.Size = 18
$r.Font.Name = "Consolas"
$r.Font.Color.RGB = 0x0055007F                 # 7F0055

$r = $tr9.Characters(8, 1)                     # " "
$r.Font.Bold = $true
$r.Font.Size = 18
$r.Font.Name = "Consolas"
$r.Font.Color.RGB = 0x00000000                 # 000000

$r = $tr9.Characters(9, 22)                    # "XSCUGIC_INTR_PRIO_MASK"
$r.Font.Bold = $true

$r = $tr9.Characters(31, 1)                    # tab
$r.Font.Bold = $true

$r = $tr9.Characters(32, 11)                   # "0x000000F8U"
$r.Font.Bold = $true

# --- 2b. New textbox "文字方塊 12" -----------------------------------
$tb12 = $slide2.Shapes.AddTextbox(1, 12.3367, 467.8593, 740.62496, 50.8922)
$tb12.Name = "文字方塊 12"
$tb12.Fill.Visible = 0
$tb12.TextFrame.WordWrap = -1
$tb12.TextFrame.AutoSize = 1

$tr12 = $tb12.TextFrame.TextRange
$line1 = "#define XSCUGIC_PENDING_SET_OFFSET`t0x00000200U /**< Pending Set Register */"
$line2 = "#define  XSCUGIC_PENDING_CLR_OFFSET`t0x00000280U /**< Pending Clear Register */"
$tr12.Text = $line1 + "`r" + $line2

$r = $tr12.Characters(1, 7)                    # "#define"
$r.Font.Bold = $true
$r.Font.Size = 18
$r.Font.Name = "Consolas"
$r.Font.Color.RGB = 0x0055007F                 # 7F0055

$r = $tr12.Characters(8, 1)                    # " "
$r.Font.Bold = $true
$r.Font.Size = 18
$r.Font.Name = "Consolas"
$r.Font.Color.RGB = 0x00000000                 # 000000

# remainder of paragraph 1 keeps default formatting (no explicit run needed)

$p2start = $line1.Length + 1 + 1               # start index of paragraph 2 (1-based)
$r = $tr12.Characters($p2start, 7)             # "#define"
$r.Font.Bold = $true
$r.Font.Size = 18
$r.Font.Name = "Consolas"
$r.Font.Color.RGB = 0x0055007F                 # 7F0055

# remainder of paragraph 2 keeps default formatting, but stays split into the
# same two runs as the source deck (language-boundary split, no visual change).
# Re-assigning identical text to a Characters() sub-range forces a run
# boundary there without touching any Font property (so no stray attribute
# such as b="0" gets written).
$tailStart = $p2start + 7
$tail = "  XSCUGIC_PENDING_CLR_OFFSET`t0x00000280U /**< Pending Clear"
$r = $tr12.Characters($tailStart, $tail.Length)
$r.Text = $tail
$tailStart2 = $tailStart + $tail.Length
$tail2 = " Register */"
$r = $tr12.Characters($tailStart2, $tail2.Length)
$r.Text = $tail2

# --- 2c. New textbox "文字方塊 14" -----------------------------------
$tb14 = $slide2.Shapes.AddTextbox(1, 21.09945, 405.3452, 557.24992, 60.5859)
$tb14.Name = "文字方塊 14"
$tb14.Fill.Visible = 0
$tb14.TextFrame.WordWrap = -1
$tb14.TextFrame.AutoSize = 1

$tr14 = $tb14.TextFrame.TextRange
$tr14.Text = "pending暫存器的addresses"
$tr14.Font.Size = 44

# force the three language-boundary runs to stay distinct (mirrors source deck)
$r = $tr14.Characters(1, 7)                    # "pending"
$r.Font.Size = 44
$r = $tr14.Characters(8, 4)                    # "暫存器的"
$r.Font.Size = 44
$r = $tr14.Characters(12, 9)                   # "addresses"
$r.Font.Size = 44

# ====================================================================
# Slide 5 — clear title & body text
# ====================================================================
$slide5 = $p.Slides.Item(5)
$slide5.Shapes.Item(1).TextFrame.TextRange.Text = ""
$slide5.Shapes.Item(2).TextFrame.TextRange.Text = ""
